$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.344
$ws.Range("B4").Value = 4.808999999999999
$ws.Range("A6").Value = -22.178
$ws.Range("A7").Value = -20.637
$ws.Range("B9").Value = 5.808
$ws.Range("B12").Value = 4.939
$ws.Range("A16").Value = -21.79
$ws.Range("B17").Value = 5.179
$ws.Range("B18").Value = 5.988999999999999
$ws.Range("B19").Value = 7.491
$ws.Range("A20").Value = -21.902
$ws.Range("B20").Value = 5.173
$ws.Range("B26").Value = 6.291
$ws.Range("A28").Value = -21.624
$ws.Range("A29").Value = -21.5
$ws.Range("B31").Value = 6.149999999999999
$ws.Range("A32").Value = -21.418
$ws.Range("B39").Value = 7.145
$ws.Range("A40").Value = -20.447
$ws.Range("B40").Value = 7.25
$ws.Range("B41").Value = 6.251
$ws.Range("B42").Value = 6.101
$ws.Range("B43").Value = 6.071000000000001
$ws.Range("A46").Value = -21.489
$ws.Range("B47").Value = 5.82
$ws.Range("B48").Value = 5.443
$ws.Range("A51").Value = -21.557
$ws.Range("A52").Value = -21.646
$ws.Range("A57").Value = -21.889
$ws.Range("A59").Value = -22.257
$ws.Range("A62").Value = -21.779
$ws.Range("B63").Value = 5.252
$ws.Range("B64").Value = 5.527
$ws.Range("A66").Value = -21.563
$ws.Range("A73").Value = -21.099
$ws.Range("A74").Value = -20.458
$ws.Range("B76").Value = 5.792
$ws.Range("B81").Value = 5.624
$ws.Range("B89").Value = 5.482
$ws.Range("A92").Value = -21.566
$ws.Range("B94").Value = 5.795
$ws.Range("A100").Value = -21.961